$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A53").Value = "ZYW2LG"
$ws.Range("B53").Value = "Kit de engranaje de acople de fusor para impresora HP"
$ws.Range("C53").Value = "LaserJet 5100/5200"
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 180000
$ws.Range("F53").Value = 2
$ws.Range("G53").Value = 0
$ws.Range("H53").Formula = "=(E53-D53)*G53"
$ws.Range("I53").Formula = "=D53*F53"
$ws.Range("J53").Value = 0
